$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F(x) calculation values in column G
$ws.Range("G2").Value = 0.047
$ws.Range("G3").Value = 0.374
$ws.Range("G4").Value = 2.543
$ws.Range("G5").Value = 7.114
$ws.Range("G6").Value = 182.349

# Update the active selection/cell to G6
$ws.Range("G6").Select()
